# Scheduled Sheets data-refresh: updates cached market-board price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) for specific leve
# rows across all class tables. Some rows switch which of the NQ/HQ profit
# columns (M/N) is populated, so those cells are explicitly cleared/created.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 296.08334
$ws.Range("I6").Value = 295.4
$ws.Range("K6").Value = 886.1999999999999
$ws.Range("M6").Value = -774.1999999999999
$ws.Range("H17").Value = 1139.7142
$ws.Range("J17").Value = 1139.7142
$ws.Range("L17").Value = 3419.1426
$ws.Range("N17").Value = -3755.1426
$ws.Range("H19").Value = 860.6667
$ws.Range("J19").Value = 613.1429000000001
$ws.Range("L19").Value = 613.1429000000001
$ws.Range("N19").Value = -963.1429000000001
$ws.Range("H31").Value = 365
$ws.Range("I31").Value = 365
$ws.Range("K31").Value = 1095
$ws.Range("M31").Value = -865
$ws.Range("H34").Value = 4007.75
$ws.Range("I34").Value = 4007.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4007.75
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3804.75
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 4007.75
$ws.Range("I36").Value = 4007.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4007.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3292.75
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 4824.081
$ws.Range("I40").Value = 3948.1035
$ws.Range("K40").Value = 3948.1035
$ws.Range("M40").Value = -3773.1035
$ws.Range("H41").Value = 2118.75
$ws.Range("I41").Value = 406.4
$ws.Range("J41").Value = 3341.8572
$ws.Range("K41").Value = 406.4
$ws.Range("L41").Value = 3341.8572
$ws.Range("M41").Value = 33.60000000000002
$ws.Range("N41").Value = -4221.8572
$ws.Range("H48").Value = 3906.6667
$ws.Range("J48").Value = 3903.6
$ws.Range("L48").Value = 11710.8
$ws.Range("N48").Value = -12294.8
$ws.Range("H56").Value = 3906.6667
$ws.Range("J56").Value = 3903.6
$ws.Range("L56").Value = 11710.8
$ws.Range("N56").Value = -12778.8
$ws.Range("H61").Value = 12300.667
$ws.Range("I61").Value = 17942.5
$ws.Range("K61").Value = 53827.5
$ws.Range("M61").Value = -53655.5
$ws.Range("H76").Value = 2341.3333
$ws.Range("J76").Value = 1809.8
$ws.Range("L76").Value = 1809.8
$ws.Range("N76").Value = -2439.8
$ws.Range("H79").Value = 2341.3333
$ws.Range("J79").Value = 1809.8
$ws.Range("L79").Value = 1809.8
$ws.Range("N79").Value = -3993.8
$ws.Range("H98").Value = 814.5
$ws.Range("I98").Value = 774.25
$ws.Range("K98").Value = 774.25
$ws.Range("M98").Value = 723.75
$ws.Range("H107").Value = 593.3333
$ws.Range("I107").Value = 560.06665
$ws.Range("K107").Value = 560.06665
$ws.Range("M107").Value = 1359.93335
$ws.Range("H122").Value = 814.5
$ws.Range("I122").Value = 774.25
$ws.Range("K122").Value = 2322.75
$ws.Range("M122").Value = 127.25
$ws.Range("H131").Value = 1012
$ws.Range("I131").Value = 1000
$ws.Range("K131").Value = 3000
$ws.Range("M131").Value = 2040
$ws.Range("H132").Value = 9392.241
$ws.Range("I132").Value = 8899
$ws.Range("K132").Value = 26697
$ws.Range("M132").Value = -24167
$ws.Range("H137").Value = 4695.1304
$ws.Range("I137").Value = 1998.9
$ws.Range("K137").Value = 5996.700000000001
$ws.Range("M137").Value = -3446.700000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 5629.7896
$ws.Range("J61").Value = 6772.1816
$ws.Range("L61").Value = 6772.1816
$ws.Range("N61").Value = -7196.1816
$ws.Range("H136").Value = 5629.7896
$ws.Range("J136").Value = 6772.1816
$ws.Range("L136").Value = 20316.5448
$ws.Range("N136").Value = -25416.5448

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5509.875
$ws.Range("I86").Value = 4868.4287
$ws.Range("K86").Value = 4868.4287
$ws.Range("M86").Value = -3745.4287
$ws.Range("H89").Value = 5509.875
$ws.Range("I89").Value = 4868.4287
$ws.Range("K89").Value = 24342.1435
$ws.Range("M89").Value = -18726.1435
$ws.Range("H100").Value = 25641
$ws.Range("J100").Value = 25641
$ws.Range("L100").Value = 25641
$ws.Range("N100").Value = -27805

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 445.16666
$ws.Range("I7").Value = 89.46154
$ws.Range("J7").Value = 1370
$ws.Range("K7").Value = 89.46154
$ws.Range("L7").Value = 1370
$ws.Range("M7").Value = 23.53846
$ws.Range("N7").Value = -1596
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4826
$ws.Range("H22").Value = 243.33333
$ws.Range("I22").Value = 243.33333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 243.33333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 106.66667
$ws.Range("N22").ClearContents()
$ws.Range("H92").Value = 9829.857
$ws.Range("J92").Value = 9829.857
$ws.Range("L92").Value = 9829.857
$ws.Range("N92").Value = -14821.857

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 437.5
$ws.Range("I107").Value = 345
$ws.Range("J107").Value = 483.75
$ws.Range("K107").Value = 1035
$ws.Range("L107").Value = 1451.25
$ws.Range("M107").Value = 885
$ws.Range("N107").Value = -5291.25
$ws.Range("H113").Value = 801.7
$ws.Range("J113").Value = 899.1429000000001
$ws.Range("L113").Value = 2697.4287
$ws.Range("N113").Value = -7037.4287
$ws.Range("H134").Value = 4646
$ws.Range("I134").Value = 1076.6666
$ws.Range("K134").Value = 3229.9998
$ws.Range("M134").Value = 1840.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H101").Value = 45995
$ws.Range("J101").Value = 45995
$ws.Range("L101").Value = 45995
$ws.Range("N101").Value = -52485
$ws.Range("H134").Value = 100000.5
$ws.Range("J134").Value = 100000.5
$ws.Range("L134").Value = 300001.5
$ws.Range("N134").Value = -305071.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 553.1667
$ws.Range("I22").Value = 522.5
$ws.Range("K22").Value = 522.5
$ws.Range("M22").Value = -227.5
$ws.Range("H27").Value = 553.1667
$ws.Range("I27").Value = 522.5
$ws.Range("K27").Value = 522.5
$ws.Range("M27").Value = -415.5
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 5998.5713
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3812
$ws.Range("H132").Value = 2597.6
$ws.Range("I132").Value = 2597.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7792.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5262.799999999999
$ws.Range("N132").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 8000
$ws.Range("I70").Value = 8000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7685
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 8000
$ws.Range("I73").Value = 8000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -6908
$ws.Range("N73").ClearContents()
$ws.Range("H136").Value = 2713.72
$ws.Range("I136").Value = 1673.1177
$ws.Range("K136").Value = 5019.3531
$ws.Range("M136").Value = -2469.3531

